$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Spiderman"
$ws.Range("B4").Value = "19513fdc9da4fb72a4a05eb66917548d3c90ff94d5419e1f2363eea89dfee1dd"
$ws.Range("C4").Value = "Fortnite@gmail.com"
$ws.Range("D4").Value = "member"

$ws.Range("A5").Value = "JohnDoe"
$ws.Range("B5").Value = "e7cf3ef4f17c3999a94f2c6f612e8a888e5b1026878e4e19398b23bd38ec221a"
$ws.Range("C5").Value = "libralinkcpsc@gmail.com"
$ws.Range("D5").Value = "member"
